$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("testresults_ISR")

# Excel alignment constants
$xlCenter = -4108
$xlLeft = -4131
$xlRight = -4152

# Helper cell used to stage numeric values before moving them into columns G/H
# whose column-level style uses a text ("@") number format; writing a number
# straight into such a cell would otherwise be coerced to a text string.
$helper = $ws.Cells.Item(1, 11)

function Set-NumericCell($row, $col, $value) {
    $helper.Value = $value
    $helper.Cut($ws.Cells.Item($row, $col)) | Out-Null
}

# ---------------------------------------------------------------------------
# Row 26 : Area1 / isr per 100000 / Byars  (new right-aligned style)
# ---------------------------------------------------------------------------
$ws.Cells.Item(26, 1).Value = "Area1"

$ws.Cells.Item(26, 2).Value = 895
$ws.Cells.Item(26, 2).HorizontalAlignment = $xlRight

$ws.Cells.Item(26, 3).Value = 840.93969634711902
$ws.Cells.Item(26, 3).HorizontalAlignment = $xlRight

$ws.Cells.Item(26, 4).Value = 9718.074126844951
$ws.Cells.Item(26, 4).HorizontalAlignment = $xlRight

$ws.Cells.Item(26, 5).Value = 10342.806245569416
$ws.Cells.Item(26, 5).HorizontalAlignment = $xlRight

$ws.Cells.Item(26, 6).Value = 9676.2084578011254
$ws.Cells.Item(26, 6).HorizontalAlignment = $xlRight

Set-NumericCell 26 7 11043.22699510916
$ws.Cells.Item(26, 7).HorizontalAlignment = $xlRight

$ws.Cells.Item(26, 8).Value = "95%"

$ws.Cells.Item(26, 9).Value = "isr per 100000"

$ws.Cells.Item(26, 10).Value = "Byars"
$ws.Cells.Item(26, 10).HorizontalAlignment = $xlLeft

# ---------------------------------------------------------------------------
# Row 27 : Area2 / isr per 100000 / Byars
# ---------------------------------------------------------------------------
$ws.Cells.Item(27, 1).Value = "Area2"

$ws.Cells.Item(27, 2).Value = 91

$ws.Cells.Item(27, 3).Value = 868.50151748646806

$ws.Cells.Item(27, 4).Value = 9718.074126844951

$ws.Cells.Item(27, 5).Value = 1018.2420269135216

$ws.Cells.Item(27, 6).Value = 819.80264950202934

Set-NumericCell 27 7 1250.1887682774923
$ws.Cells.Item(27, 7).HorizontalAlignment = $xlCenter
$ws.Cells.Item(27, 7).NumberFormat = "@"

$ws.Cells.Item(27, 8).Value = "95%"

$ws.Cells.Item(27, 9).Value = "isr per 100000"

$ws.Cells.Item(27, 10).Value = "Byars"

# ---------------------------------------------------------------------------
# Row 28 : Area1 / smr x 100 / Byars  (no age specific denominator -> D blank)
# ---------------------------------------------------------------------------
$ws.Cells.Item(28, 1).Value = "Area1"

$ws.Cells.Item(28, 2).Value = 895

$ws.Cells.Item(28, 3).Value = 840.93969634711902

$ws.Cells.Item(28, 5).Value = 106.42855889521074

$ws.Cells.Item(28, 6).Value = 99.569197883270135

Set-NumericCell 28 7 113.63596172418198
$ws.Cells.Item(28, 7).HorizontalAlignment = $xlCenter
$ws.Cells.Item(28, 7).NumberFormat = "@"

$ws.Cells.Item(28, 8).Value = "95%"

$ws.Cells.Item(28, 9).Value = "smr x 100"

$ws.Cells.Item(28, 10).Value = "Byars"

# ---------------------------------------------------------------------------
# Row 29 : Area2 / smr x 100 / Byars  (no age specific denominator -> D blank)
# ---------------------------------------------------------------------------
$ws.Cells.Item(29, 1).Value = "Area2"

$ws.Cells.Item(29, 2).Value = 91

$ws.Cells.Item(29, 3).Value = 868.50151748646806

$ws.Cells.Item(29, 5).Value = 10.477817040937738

$ws.Cells.Item(29, 6).Value = 8.4358550758264759

Set-NumericCell 29 7 12.864573288487314
$ws.Cells.Item(29, 7).HorizontalAlignment = $xlCenter
$ws.Cells.Item(29, 7).NumberFormat = "@"

$ws.Cells.Item(29, 8).Value = "95%"

$ws.Cells.Item(29, 9).Value = "smr x 100"

$ws.Cells.Item(29, 10).Value = "Byars"

# ---------------------------------------------------------------------------
# Update the sheet view selection to match the new active cell
# ---------------------------------------------------------------------------
$ws.Range("A27").Select() | Out-Null
